# Hortaliza, Terminal La Palmera de La Serena - Coliflor
# A new weekly record (fecha 2022-07-11 / serial 44753) is inserted at the
# top of the data block (row 709), pushing all the existing records for
# this market down by two rows (one row each for "Primera" and "Segunda"
# quality).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current first data row of this
# block (row 709); this shifts rows 709:744 down to 711:746.
$ws.Rows("709:710").Insert()

# --- New row 709 ("Primera") ---
$ws.Cells.Item(709, 1).Value = 8
$ws.Cells.Item(709, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(709, 3).Value = "Coquimbo"
$ws.Cells.Item(709, 4).Value = 44753
$ws.Cells.Item(709, 5).Value = 4
$ws.Cells.Item(709, 6).Value = 100112008
$ws.Cells.Item(709, 7).Value = "Coliflor"
$ws.Cells.Item(709, 8).Value = "Sin especificar"
$ws.Cells.Item(709, 9).Value = "Primera"
$ws.Cells.Item(709, 10).Value = 2400
$ws.Cells.Item(709, 11).Value = 800
$ws.Cells.Item(709, 12).Value = 900
$ws.Cells.Item(709, 13).Value = 850
$ws.Cells.Item(709, 14).Value = "`$/unidad"
$ws.Cells.Item(709, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(709, 16).Value = 850
$ws.Cells.Item(709, 17).Value = 1
$ws.Cells.Item(709, 18).Value = "Hortaliza"

# --- New row 710 ("Segunda") ---
$ws.Cells.Item(710, 1).Value = 8
$ws.Cells.Item(710, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(710, 3).Value = "Coquimbo"
$ws.Cells.Item(710, 4).Value = 44753
$ws.Cells.Item(710, 5).Value = 4
$ws.Cells.Item(710, 6).Value = 100112008
$ws.Cells.Item(710, 7).Value = "Coliflor"
$ws.Cells.Item(710, 8).Value = "Sin especificar"
$ws.Cells.Item(710, 9).Value = "Segunda"
$ws.Cells.Item(710, 10).Value = 1360
$ws.Cells.Item(710, 11).Value = 700
$ws.Cells.Item(710, 12).Value = 750
$ws.Cells.Item(710, 13).Value = 725
$ws.Cells.Item(710, 14).Value = "`$/unidad"
$ws.Cells.Item(710, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(710, 16).Value = 725
$ws.Cells.Item(710, 17).Value = 1
$ws.Cells.Item(710, 18).Value = "Hortaliza"
